$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header columns (row 1) to short machine-friendly field names
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'

# 2. Title-case the Spanish connector words (de/del/el/la/las/los/y) in
#    the municipality/state names across columns A and B
$nameUpdates = @(
    @('B22', 'Amatenango De La Frontera'),
    @('B24', 'Bejucal De Ocampo'),
    @('B30', 'Chiapa De Corzo'),
    @('B33', 'Comitán De Domínguez'),
    @('B49', 'Mazapa De Madero'),
    @('B52', 'Montecristo De Guerrero'),
    @('B55', 'Ocozocoautla De Espinosa'),
    @('B61', 'Salto De Agua'),
    @('B62', 'San Cristóbal De Las Casas'),
    @('B93', 'San Juan De Sabinas'),
    @('A102', 'Ciudad De México'),
    @('B106', 'Cuajimalpa De Morelos'),
    @('A128', 'Estado De México'),
    @('B128', 'Acambay De Ruíz Castañeda'),
    @('B131', 'Almoloya De Alquisiras'),
    @('B132', 'Almoloya De Juárez'),
    @('B137', 'Atizapán De Zaragoza'),
    @('B142', 'Chapa De Mota'),
    @('B145', 'Coacalco De Berriozábal'),
    @('B150', 'Ecatepec De Morelos'),
    @('B154', 'Ixtapan De La Sal'),
    @('B155', 'Ixtapan Del Oro'),
    @('B166', 'Naucalpan De Juárez'),
    @('B172', 'San Felipe Del Progreso'),
    @('B174', 'San Simón De Guerrero'),
    @('B182', 'Tenango Del Valle'),
    @('B190', 'Tlalnepantla De Baz'),
    @('B194', 'Valle De Bravo'),
    @('B195', 'Valle De Chalco Solidaridad'),
    @('B196', 'Villa De Allende'),
    @('B208', 'Apaseo El Alto'),
    @('B209', 'Apaseo El Grande'),
    @('B214', 'Dolores Hidalgo Cuna De La Independencia Nacional'),
    @('B225', 'San Diego De La Unión'),
    @('B227', 'San Francisco Del Rincón'),
    @('B228', 'San Luis De La Paz'),
    @('B229', 'Santa Cruz De Juventino Rosas'),
    @('B231', 'Silao De La Victoria'),
    @('B235', 'Valle De Santiago'),
    @('B240', 'Acapulco De Juárez'),
    @('B242', 'Ajuchitlán Del Progreso'),
    @('B243', 'Alcozauca De Guerrero'),
    @('B247', 'Atenango Del Río'),
    @('B248', 'Atlamajalcingo Del Monte'),
    @('B250', 'Atoyac De Álvarez'),
    @('B251', 'Ayutla De Los Libres'),
    @('B253', 'Buenavista De Cuéllar'),
    @('B254', 'Chilapa De Álvarez'),
    @('B255', 'Chilpancingo De Los Bravo'),
    @('B256', 'Coahuayutla De José María Izazaga'),
    @('B261', 'Coyuca De Benítez'),
    @('B262', 'Coyuca De Catalán'),
    @('B266', 'Cuetzala Del Progreso'),
    @('B267', 'Cutzamala De Pinzón'),
    @('B273', 'Huitzuco De Los Figueroa'),
    @('B274', 'Iguala De La Independencia'),
    @('B276', 'Ixcateopan De Cuauhtémoc'),
    @('B277', 'Zihuatanejo De Azueta'),
    @('B279', 'La Unión De Isidoro Montes De Oca'),
    @('B282', 'Mártir De Cuilapan'),
    @('B295', 'Taxco De Alarcón'),
    @('B297', 'Técpan De Galeana'),
    @('B299', 'Tepecoacuilco De Trujano'),
    @('B301', 'Tixtla De Guerrero'),
    @('B305', 'Tlalixtaquilla De Maldonado'),
    @('B306', 'Tlapa De Comonfort'),
    @('B317', 'Agua Blanca De Iturbide'),
    @('B321', 'Atotonilco De Tula'),
    @('B322', 'Atotonilco El Grande'),
    @('B327', 'Cuautepec De Hinojosa'),
    @('B330', 'Huasca De Ocampo'),
    @('B334', 'Huejutla De Reyes'),
    @('B337', 'Jacala De Ledezma'),
    @('B343', 'Mineral Del Chico'),
    @('B344', 'Mixquiahuala De Juárez'),
    @('B345', 'Nopala De Villagrán'),
    @('B346', 'Omitlán De Juárez'),
    @('B347', 'Pachuca De Soto'),
    @('B350', 'Progreso De Obregón'),
    @('B353', 'Santiago De Anaya'),
    @('B354', 'Santiago Tulantepec De Lugo Guerrero'),
    @('B358', 'Tenango De Doria'),
    @('B359', 'Tepehuacán De Guerrero'),
    @('B362', 'Tezontepec De Aldama'),
    @('B367', 'Tula De Allende'),
    @('B368', 'Tulancingo De Bravo'),
    @('B371', 'Zacualtipán De Ángeles'),
    @('B381', 'Encarnación De Díaz'),
    @('B385', 'Lagos De Moreno'),
    @('B390', 'San Martín De Bolaños'),
    @('B392', 'San Miguel El Alto'),
    @('B394', 'Tamazula De Gordiano'),
    @('B398', 'Tepatitlán De Morelos'),
    @('B402', 'Valle De Guadalupe'),
    @('B403', 'Zacoalco De Torres'),
    @('B406', 'Zapotlán El Grande'),
    @('B421', 'Coalcomán De Vázquez Pallares'),
    @('B463', 'Tiquicheo De Nicolás Romero'),
    @('B481', 'Coatlán Del Río'),
    @('B488', 'Jonacatepec De Leandro Valle'),
    @('B491', 'Puente De Ixtla'),
    @('B496', 'Tetela Del Volcán'),
    @('B498', 'Tlaltizapán De Zapata'),
    @('B506', 'Zacualpan De Amilpas'),
    @('B521', 'Acatlán De Pérez Figueroa'),
    @('B528', 'Ayoquezco De Aldama'),
    @('B531', 'Capulálpam De Méndez'),
    @('B533', 'Chalcatongo De Hidalgo'),
    @('B534', 'Ciénega De Zimatlán'),
    @('B536', 'Coicoyán De Las Flores'),
    @('B539', 'Constancia Del Rosario'),
    @('B542', 'Cuilápam De Guerrero'),
    @('B544', 'Eloxochitlán De Flores Magón'),
    @('B545', 'Fresnillo De Trujano'),
    @('B546', 'Guadalupe De Ramírez'),
    @('B548', 'Heroica Ciudad De Ejutla De Crespo'),
    @('B549', 'Heroica Ciudad De Huajuapan De León'),
    @('B550', 'Heroica Ciudad De Tlaxiaco'),
    @('B551', 'Huautla De Jiménez'),
    @('B553', 'Ixtlán De Juárez'),
    @('B554', 'Heroica Ciudad De Juchitán De Zaragoza'),
    @('B563', 'Mariscala De Juárez'),
    @('B564', 'Mártires De Tacubaya'),
    @('B567', 'Miahuatlán De Porfirio Díaz'),
    @('B568', 'Mixistlán De La Reforma'),
    @('B572', 'Nejapa De Madero'),
    @('B574', 'Oaxaca De Juárez'),
    @('B575', 'Ocotlán De Morelos'),
    @('B576', 'Pinotepa De Don Luis'),
    @('B578', 'Putla Villa De Guerrero'),
    @('B579', 'Reforma De Pineda'),
    @('B584', 'San Agustín De Las Juntas'),
    @('B594', 'San Antonino El Alto'),
    @('B600', 'San Baltazar Yatzachi El Bajo'),
    @('B605', 'San Felipe Jalapa De Díaz'),
    @('B623', 'San José Del Progreso'),
    @('B630', 'San Juan Bautista Lo De Soto'),
    @('B638', 'San Juan Del Estado'),
    @('B639', 'San Juan Del Río'),
    @('B673', 'San Mateo Del Mar'),
    @('B683', 'San Miguel Del Puerto'),
    @('B684', 'San Miguel Del Río'),
    @('B686', 'San Miguel El Grande'),
    @('B697', 'San Pablo Villa De Mitla'),
    @('B701', 'San Pedro El Alto'),
    @('B731', 'Santa Ana Del Valle'),
    @('B746', 'Santa Cruz Tacache De Mina'),
    @('B766', 'Santa María Jalapa Del Marqués'),
    @('B789', 'Santiago Del Río'),
    @('B818', 'Santo Domingo De Morelos'),
    @('B839', 'Tamazulápam Del Espíritu Santo'),
    @('B840', 'Tanetze De Zaragoza'),
    @('B841', 'Tataltepec De Valdés'),
    @('B842', 'Teococuilco De Marcos Pérez'),
    @('B843', 'Teotitlán De Flores Magón'),
    @('B844', 'Teotitlán Del Valle'),
    @('B845', 'Heroica Villa Tezoatlán De Segura Y Luna, Cuna De La Independencia De Oaxaca'),
    @('B846', 'Tlacolula De Matamoros'),
    @('B847', 'Tlalixtac De Cabrera'),
    @('B848', 'Totontepec Villa De Morelos'),
    @('B851', 'Villa De Chilapa De Díaz'),
    @('B852', 'Villa De Etla'),
    @('B853', 'Villa De Tututepec'),
    @('B854', 'Villa De Zaachila'),
    @('B857', 'Villa Sola De Vega'),
    @('B858', 'Villa Talea De Castro'),
    @('B859', 'Zapotitlán Del Río'),
    @('B861', 'Zimatlán De Álvarez'),
    @('B886', 'Ayotoxco De Guerrero'),
    @('B898', 'Chila De La Sal'),
    @('B909', 'Cuapiaxtla De Madero'),
    @('B912', 'Cuayuca De Andrade'),
    @('B913', 'Cuetzalan Del Progreso'),
    @('B928', 'Huehuetlán El Chico'),
    @('B929', 'Huehuetlán El Grande'),
    @('B933', 'Huitzilan De Serdán'),
    @('B935', 'Ixcamilpa De Guerrero'),
    @('B938', 'Izúcar De Matamoros'),
    @('B949', 'Los Reyes De Juárez'),
    @('B950', 'Mazapiltepec De Juárez'),
    @('B962', 'Palmar De Bravo'),
    @('B972', 'San Diego La Mesa Tochimiltzingo'),
    @('B983', 'San Nicolás De Los Ranchos'),
    @('B987', 'San Salvador El Seco'),
    @('B988', 'San Salvador El Verde'),
    @('B994', 'Tecali De Herrera'),
    @('B1002', 'Tepanco De López'),
    @('B1003', 'Tepango De Rodríguez'),
    @('B1004', 'Tepatlaxco De Hidalgo'),
    @('B1010', 'Tepexi De Rodríguez'),
    @('B1012', 'Tetela De Ocampo'),
    @('B1013', 'Teteles De Avila Castillo'),
    @('B1018', 'Tlacotepec De Benito Juárez'),
    @('B1029', 'Totoltepec De Guerrero'),
    @('B1031', 'Tuzamapan De Galeana'),
    @('B1035', 'Xayacatlán De Bravo'),
    @('B1041', 'Xochitlán De Vicente Suárez'),
    @('B1055', 'Amealco De Bonfil'),
    @('B1057', 'Cadereyta De Montes'),
    @('B1058', 'Jalpan De Serra'),
    @('B1059', 'Landa De Matamoros'),
    @('B1060', 'Pinal De Amoles'),
    @('B1062', 'San Juan Del Río'),
    @('B1070', 'Axtla De Terrazas'),
    @('B1077', 'Mexquitic De Carmona'),
    @('B1086', 'Villa De Arista'),
    @('B1087', 'Villa De Ramos'),
    @('B1110', 'Jalpa De Méndez'),
    @('B1128', 'Acuamanala De Miguel Hidalgo'),
    @('B1134', 'Contla De Juan Cuamatzi'),
    @('B1142', 'Ixtacuixtla De Mariano Matamoros'),
    @('B1145', 'Mazatecochco De José María Morelos'),
    @('B1148', 'Papalotla De Xicohténcatl'),
    @('B1152', 'San Pablo Del Monte'),
    @('B1157', 'Tepetitla De Lardizábal'),
    @('B1160', 'Tetla De La Solidaridad'),
    @('B1177', 'Alto Lucero De Gutiérrez Barrios'),
    @('B1181', 'Amatlán De Los Reyes'),
    @('B1191', 'Castillo De Teayo'),
    @('B1193', 'Cazones De Herrera'),
    @('B1202', 'Cosamaloapan De Carpio'),
    @('B1215', 'Hueyapan De Ocampo'),
    @('B1216', 'Ignacio De La Llave'),
    @('B1219', 'Ixhuatlán De Madero'),
    @('B1220', 'Ixhuatlán Del Café'),
    @('B1230', 'Juchique De Ferrer'),
    @('B1234', 'Las Vigas De Ramírez'),
    @('B1235', 'Lerdo De Tejada'),
    @('B1239', 'Martínez De La Torre'),
    @('B1241', 'Medellín De Bravo'),
    @('B1253', 'Paso De Ovejas'),
    @('B1256', 'Poza Rica De Hidalgo'),
    @('B1262', 'Sayula De Alemán'),
    @('B1267', 'Tatahuicapan De Juárez'),
    @('B1288', 'Vega De Alatorre'),
    @('B1294', 'Zontecomatlán De López Y Fuentes'),
    @('B1306', 'Nochistlán De Mejía')
)
foreach ($pair in $nameUpdates) {
    $ws.Range($pair[0]).Value = $pair[1]
}

# 3. Remove the trailing footer/metadata rows (1318-1322) that are no
#    longer part of the cleaned dataset
$ws.Range("A1318:A1322").EntireRow.Delete()

